# BestBall DK5.xlsx update
# - Append a new mock-draft block (draft #13, rows 242-261) to Sheet1
# - Re-point the sheet view's selection to the new last cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Seed formatting for the new rows by copying the previous draft
#    block's date column so the new I-column cells pick up the existing
#    short-date style instead of minting a new number format.
# ---------------------------------------------------------------------
$ws.Range("I222:I241").Copy($ws.Range("I242:I261"))

# ---------------------------------------------------------------------
# 2. New draft data (Team 13 - "Anchor RB" / "Elite QB/TE")
# ---------------------------------------------------------------------
$draftRows = @(
    @{ Row=242; Pos="WR"; Name="Jamar Chase";        Team="CIN"; Pick=1;  F=6;   G=4.8 },
    @{ Row=243; Pos="TE"; Name="Mark Andrews";       Team="BAL"; Pick=2;  F=19;  G=18.5 },
    @{ Row=244; Pos="RB"; Name="Javonte Williams";   Team="DEN"; Pick=3;  F=30;  G=26.2 },
    @{ Row=245; Pos="QB"; Name="Lamar Jackson";      Team="BAL"; Pick=4;  F=43;  G=44.1 },
    @{ Row=246; Pos="WR"; Name="Chris Godwin";       Team="TB";  Pick=5;  F=54;  G=48.3 },
    @{ Row=247; Pos="WR"; Name="Rashod Bateman";     Team="BAL"; Pick=6;  F=67;  G=68.6 },
    @{ Row=248; Pos="WR"; Name="Hunter Renfrow";     Team="LV";  Pick=7;  F=78;  G=79.1 },
    @{ Row=249; Pos="WR"; Name="Brandon Aiyuk";      Team="SF";  Pick=8;  F=91;  G=89.3 },
    @{ Row=250; Pos="QB"; Name="Trey Lance";         Team="SF";  Pick=9;  F=102; G=98.1 },
    @{ Row=251; Pos="WR"; Name="Chase Claypool";     Team="PIT"; Pick=10; F=115; G=111 },
    @{ Row=252; Pos="TE"; Name="Pat Freiermuth";     Team="PIT"; Pick=11; F=126; G=118.2 },
    @{ Row=253; Pos="RB"; Name="Isaiah Spiller";     Team="LAC"; Pick=12; F=139; G=134.3 },
    @{ Row=254; Pos="WR"; Name="Jamison Crowder";    Team="BUF"; Pick=13; F=150; G=148.1 },
    @{ Row=255; Pos="RB"; Name="Darrell Henderson";  Team="LAR"; Pick=14; F=163; G=151.9 },
    @{ Row=256; Pos="RB"; Name="Khalil Herbert";     Team="CHI"; Pick=15; F=174; G=172 },
    @{ Row=257; Pos="RB"; Name="Tyrion Davis-Price"; Team="SF";  Pick=16; F=187; G=177 },
    @{ Row=258; Pos="RB"; Name="Brian Robinson";     Team="WAS"; Pick=17; F=198; G=191.9 },
    @{ Row=259; Pos="RB"; Name="James White";        Team="NE";  Pick=18; F=211; G=211 },
    @{ Row=260; Pos="WR"; Name="Danny Gray";         Team="SF";  Pick=19; F=222; G=222 },
    @{ Row=261; Pos="RB"; Name="Sony Michel";        Team="MIA"; Pick=20; F=235; G=235 }
)

foreach ($r in $draftRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 13
    $ws.Cells.Item($row, 2).Value = $r.Pos
    $ws.Cells.Item($row, 3).Value = $r.Name
    $ws.Cells.Item($row, 4).Value = $r.Team
    $ws.Cells.Item($row, 5).Value = $r.Pick
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Formula = "=F" + $row + "-G" + $row
    $ws.Cells.Item($row, 9).Value = 44697
    $ws.Cells.Item($row, 10).Formula = "=SUM(H`$242:H`$261)"
    $ws.Cells.Item($row, 11).Value = "Anchor RB"
    $ws.Cells.Item($row, 12).Value = "Elite QB/TE"
}

# ---------------------------------------------------------------------
# 3. Match the saved view: scrolled near the bottom, new last cell active
# ---------------------------------------------------------------------
$ws.Range("F261").Select() | Out-Null
